$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new recipe rows below the existing data (rows 1-15 -> now 1-17)
# Shared strings are registered in first-use order, so set A17 (crème) before
# A16 (tomate) to match the target uniqueCount order: 15=crème, 16=tomate.
$ws.Range("A17").Value = "panini pizz base crème"
$ws.Range("A16").Value = "panini pizz base tomate"
